# Handback status report run: refresh the "Correspond Handoff Datetime" and
# "Correspond Handback DateTime" timestamps for the b21c5673-...-e1281b4b...
# string entry on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-24 09:48:30"
$wsZhCn.Range("H4").Value = "2016-03-24 09:48:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-24 09:48:34"
$wsDeDe.Range("H4").Value = "2016-03-24 09:49:06"
